# Updates cryptos list figures (Price / Volume(1h) columns, and two
# row swaps in the Coin/Link columns) to match the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.Value = '''36.610.69'
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = '''  +0.50%  '
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.Value = '''1.958.85'
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = '''  +0.80%  '
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("E4")
$cell.Value = '''  +0.09%  '
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("D5")
$cell.Value = '''243.56'
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = '''  +0.06%  '
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.Value = '''0.616'
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = '''  +0.43%  '
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.Value = '''60.46'
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = '''  +5.42%  '
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("D8")
$cell.Value = '''1.00'
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("D9")
$cell.Value = '''0.374'
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = '''  +3.76%  '
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("D10")
$cell.Value = '''0.0787'
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = '''  -6.73%  '
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("E11")
$cell.Value = '''  +0.17%  '
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("D12")
$cell.Value = '''14.24'
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = '''  +5.45%  '
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("D13")
$cell.Value = '''21.81'
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = '''  +2.13%  '
$cell.Style = "Normal"

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$cell = $ws.Range("D14")
$cell.Value = '''2.247.94'
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = '''  +0.86%  '
$cell.Style = "Normal"

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range("D15")
$cell.Value = '''0.829'
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = '''  +2.14%  '
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.Value = '''5.26'
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = '''  +1.91%  '
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.Value = '''1.963.93'
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = '''  +1.23%  '
$cell.Style = "Normal"

# Row 18
$cell = $ws.Range("D18")
$cell.Value = '''36.466.02'
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = '''  +0.22%  '
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.Value = '''69.58'
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = '''  +0.46%  '
$cell.Style = "Normal"

# Row 20
$cell = $ws.Range("D20")
$cell.Value = '''0.0₃0850'
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = '''  -1.41%  '
$cell.Style = "Normal"

# Row 21
$cell = $ws.Range("D21")
$cell.Value = '''229.37'
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = '''  +0.37%  '
$cell.Style = "Normal"

# Row 22
$cell = $ws.Range("D22")
$cell.Value = '''5.06'
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = '''  +0.88%  '
$cell.Style = "Normal"

# Row 23
$cell = $ws.Range("E23")
$cell.Value = '''  -0.01%  '
$cell.Style = "Normal"

# Row 24
$cell = $ws.Range("D24")
$cell.Value = '''2.43'
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = '''  +3.05%  '
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("D25")
$cell.Value = '''2.34'
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = '''  +2.20%  '
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.Value = '''0.141'
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = '''  +5.23%  '
$cell.Style = "Normal"

# Row 27
$cell = $ws.Range("D27")
$cell.Value = '''9.18'
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = '''  -0.27%  '
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("D28")
$cell.Value = '''160.97'
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = '''  -0.57%  '
$cell.Style = "Normal"

# Row 29
$cell = $ws.Range("D29")
$cell.Value = '''19.31'
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = '''  +0.61%  '
$cell.Style = "Normal"

# Row 30
$cell = $ws.Range("D30")
$cell.Value = '''1.32'
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = '''  +20.76%  '
$cell.Style = "Normal"

# Row 31
$cell = $ws.Range("E31")
$cell.Value = '''  +1.49%  '
$cell.Style = "Normal"

# Row 32
$cell = $ws.Range("D32")
$cell.Value = '''4.76'
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = '''  +3.90%  '
$cell.Style = "Normal"

# Row 33
$cell = $ws.Range("D33")
$cell.Value = '''0.0614'
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = '''  -0.57%  '
$cell.Style = "Normal"

# Row 34
$cell = $ws.Range("D34")
$cell.Value = '''4.41'
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = '''  +5.31%  '
$cell.Style = "Normal"

# Row 35
$cell = $ws.Range("D35")
$cell.Value = '''3.48'
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = '''  +9.56%  '
$cell.Style = "Normal"

# Row 36
$ws.Range("B36").Value = 'BinanceUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range("D36")
$cell.Value = '''1.00'
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = '''  +0.22%  '
$cell.Style = "Normal"

# Row 37
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range("D37")
$cell.Value = '''2.26'
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = '''  +4.16%  '
$cell.Style = "Normal"

# Row 38
$cell = $ws.Range("E38")
$cell.Value = '''  -1.21%  '
$cell.Style = "Normal"

# Row 39
$cell = $ws.Range("D39")
$cell.Value = '''5.47'
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = '''  -12.62%  '
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.Value = '''0.0975'
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = '''  -1.35%  '
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("E41")
$cell.Value = '''  +0.55%  '
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("D42")
$cell.Value = '''1.16'
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = '''  +0.97%  '
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.Value = '''0.0209'
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = '''  +0.08%  '
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.Value = '''15.89'
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = '''  -0.54%  '
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("D45")
$cell.Value = '''1.366.91'
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = '''  +1.87%  '
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("D46")
$cell.Value = '''88.67'
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = '''  +2.36%  '
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("E47")
$cell.Value = '''  +0.00%  '
$cell.Style = "Normal"

# Row 48
$cell = $ws.Range("D48")
$cell.Value = '''7.11'
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = '''  -1.32%  '
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("D49")
$cell.Value = '''2.83'
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = '''  +0.01%  '
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.Value = '''45.36'
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = '''  +4.77%  '
$cell.Style = "Normal"

# Row 51
$cell = $ws.Range("D51")
$cell.Value = '''2.136.92'
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = '''  +0.83%  '
$cell.Style = "Normal"
